$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price/volume figures, keyed by row number.
# "D" = Price column (stored as text, e.g. "67.888.39"), "E" = Volume(1h) column.
$updates = @(
    [PSCustomObject]@{ Row = 2; D = "67.888.39"; E = "  -2.27%  " }
    [PSCustomObject]@{ Row = 3; D = "3.818.85"; E = "  +1.40%  " }
    [PSCustomObject]@{ Row = 4; D = "0.999"; E = "  -0.15%  " }
    [PSCustomObject]@{ Row = 5; D = "603.08"; E = "  -1.93%  " }
    [PSCustomObject]@{ Row = 6; D = "172.53"; E = "  -3.28%  " }
    [PSCustomObject]@{ Row = 7; D = "3.818.02"; E = "  +1.38%  " }
    [PSCustomObject]@{ Row = 8; D = $null; E = "  -0.04%  " }
    [PSCustomObject]@{ Row = 9; D = $null; E = "  +0.92%  " }
    [PSCustomObject]@{ Row = 10; D = "0.160"; E = "  -4.44%  " }
    [PSCustomObject]@{ Row = 11; D = "6.21"; E = "  -6.60%  " }
    [PSCustomObject]@{ Row = 12; D = "0.468"; E = "  -3.27%  " }
    [PSCustomObject]@{ Row = 13; D = "38.87"; E = "  -3.11%  " }
    [PSCustomObject]@{ Row = 14; D = "0.0000245"; E = "  -3.33%  " }
    [PSCustomObject]@{ Row = 15; D = "4.450.80"; E = "  +1.32%  " }
    [PSCustomObject]@{ Row = 16; D = "3.814.18"; E = "  +1.25%  " }
    [PSCustomObject]@{ Row = 17; D = "67.828.29"; E = "  -2.38%  " }
    [PSCustomObject]@{ Row = 18; D = "7.26"; E = "  -3.56%  " }
    [PSCustomObject]@{ Row = 19; D = $null; E = "  -3.81%  " }
    [PSCustomObject]@{ Row = 20; D = "17.35"; E = "  +6.35%  " }
    [PSCustomObject]@{ Row = 21; D = "493.07"; E = "  -3.41%  " }
    [PSCustomObject]@{ Row = 22; D = "9.17"; E = "  -1.99%  " }
    [PSCustomObject]@{ Row = 23; D = "0.745"; E = "  +2.38%  " }
    [PSCustomObject]@{ Row = 24; D = "86.15"; E = "  -0.34%  " }
    [PSCustomObject]@{ Row = 25; D = "2.39"; E = "  -5.42%  " }
    [PSCustomObject]@{ Row = 26; D = $null; E = "  +7.06%  " }
    [PSCustomObject]@{ Row = 27; D = "12.38"; E = "  -3.39%  " }
    [PSCustomObject]@{ Row = 28; D = "10.28"; E = "  -2.97%  " }
    [PSCustomObject]@{ Row = 29; D = $null; E = "  +0.14%  " }
    [PSCustomObject]@{ Row = 30; D = "2.98"; E = "  +0.82%  " }
    [PSCustomObject]@{ Row = 31; D = "2.45"; E = "  -2.19%  " }
    [PSCustomObject]@{ Row = 32; D = "32.97"; E = "  +7.57%  " }
    [PSCustomObject]@{ Row = 33; D = "7.83"; E = "  -2.05%  " }
    [PSCustomObject]@{ Row = 34; D = $null; E = "  -3.74%  " }
    [PSCustomObject]@{ Row = 35; D = "0.997"; E = "  -0.25%  " }
    [PSCustomObject]@{ Row = 36; D = $null; E = "  -3.25%  " }
    [PSCustomObject]@{ Row = 37; D = "5.84"; E = "  -5.03%  " }
    [PSCustomObject]@{ Row = 38; D = "464.69"; E = "  +2.48%  " }
    [PSCustomObject]@{ Row = 39; D = "0.331"; E = "  -2.75%  " }
    [PSCustomObject]@{ Row = 40; D = $null; E = "  -5.27%  " }
    [PSCustomObject]@{ Row = 41; D = $null; E = "  -1.43%  " }
    [PSCustomObject]@{ Row = 42; D = "2.02"; E = "  -2.51%  " }
    [PSCustomObject]@{ Row = 43; D = "2.86"; E = "  -4.51%  " }
    [PSCustomObject]@{ Row = 44; D = "8.43"; E = "  -1.54%  " }
    [PSCustomObject]@{ Row = 45; D = "41.09"; E = "  -8.13%  " }
    [PSCustomObject]@{ Row = 46; D = $null; E = "  -0.02%  " }
    [PSCustomObject]@{ Row = 47; D = "2.851.22"; E = "  -3.54%  " }
    [PSCustomObject]@{ Row = 48; D = "139.15"; E = "  +0.15%  " }
    [PSCustomObject]@{ Row = 49; D = "0.0352"; E = "  -1.99%  " }
    [PSCustomObject]@{ Row = 50; D = "26.06"; E = "  -4.37%  " }
    [PSCustomObject]@{ Row = 51; D = "23.94"; E = "  +8.80%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # Force the Price cell to remain plain text so strings such as
        # "67.888.39" or "0.160" are not re-interpreted as numbers/dates.
        $cellD = $ws.Range("D" + $u.Row)
        $cellD.NumberFormat = "@"
        $cellD.Value = $u.D
        $cellD.Style = "Normal"
    }
    if ($null -ne $u.E) {
        $ws.Range("E" + $u.Row).Value = $u.E
    }
}
